$d = $word.ActiveDocument

# Update the date heading.
$d.Content.Find.Execute("2024-08-11 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-08-12 Monday", 2)

# Update the division problems in the table. Cells are addressed by
# (row, column) so duplicate problem text (e.g. "34÷4=") is handled
# unambiguously.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "27÷9="
$t.Cell(1,3).Range.Text = "99÷2="
$t.Cell(1,4).Range.Text = "50÷2="
$t.Cell(1,5).Range.Text = "62÷7="

$t.Cell(5,1).Range.Text = "29÷8="
$t.Cell(5,2).Range.Text = "59÷4="
$t.Cell(5,3).Range.Text = "22÷8="
$t.Cell(5,4).Range.Text = "43÷7="
$t.Cell(5,5).Range.Text = "54÷7="

$t.Cell(9,1).Range.Text = "78÷6="
$t.Cell(9,2).Range.Text = "76÷6="
$t.Cell(9,3).Range.Text = "74÷2="
$t.Cell(9,4).Range.Text = "91÷6="
$t.Cell(9,5).Range.Text = "79÷7="

$t.Cell(13,1).Range.Text = "21÷9="
$t.Cell(13,2).Range.Text = "15÷8="
$t.Cell(13,3).Range.Text = "58÷7="
$t.Cell(13,4).Range.Text = "64÷4="
$t.Cell(13,5).Range.Text = "48÷6="

$t.Cell(17,1).Range.Text = "77÷6="
$t.Cell(17,2).Range.Text = "82÷4="
$t.Cell(17,3).Range.Text = "92÷9="
$t.Cell(17,4).Range.Text = "94÷6="
$t.Cell(17,5).Range.Text = "58÷9="
